$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 43: uren bijgewerkt (0.5 -> 2.5) ---
$ws.Range("H43").Value = 2.5
$ws.Range("J43").Value = 2.5

# --- Row 44: nieuwe taak "Voortgangs bespreking en taken uitvoeren" ingevuld ---
$ws.Range("A44").Value = 5
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 1
$ws.Range("M44").Formula = "=C44+E44+G44+I44+K44"

# --- Row 48: uren bijgewerkt ---
$ws.Range("J48").Value = 1.5

# --- Row 49: uren bijgewerkt ---
$ws.Range("D49").Value = 0.5
$ws.Range("F49").Value = 0.5
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0.5

# --- Row 50: agendapunt gewijzigd naar "Donderdag bijwerken RA en SA" + uren bijgewerkt ---
$ws.Range("B50").Value = "Donderdag bijwerken RA en SA"
$ws.Range("D50").Value = 1.5
$ws.Range("F50").Value = 1.5
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 1.5

# --- Row 51: uren bijgewerkt ---
$ws.Range("D51").Value = 1
$ws.Range("F51").Value = 1
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 1

# --- Row 55: uren bijgewerkt ---
$ws.Range("C55").Value = 5
$ws.Range("E55").Value = 5
$ws.Range("G55").Value = 5
$ws.Range("I55").Value = 0

# --- Weergave: scroll/selectie verplaatst naar de onderkant van het blad ---
$ws.Range("B57").Select()
